$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 952.48486
$ws.Range("J17").Value = 1033.2858
$ws.Range("L17").Value = 3099.8574
$ws.Range("N17").Value = -3435.8574

$ws.Range("H28").Value = 271.25
$ws.Range("I28").Value = 209.33333
$ws.Range("J28").Value = 1200
$ws.Range("K28").Value = 209.33333
$ws.Range("L28").Value = 1200
$ws.Range("M28").Value = 275.66667
$ws.Range("N28").Value = -2170

$ws.Range("H86").Value = 93290.25
$ws.Range("I86").Value = 158357.58
$ws.Range("J86").Value = 2196
$ws.Range("K86").Value = 158357.58
$ws.Range("L86").Value = 2196
$ws.Range("M86").Value = -157234.58
$ws.Range("N86").Value = -4442

$ws.Range("H89").Value = 93290.25
$ws.Range("I89").Value = 158357.58
$ws.Range("J89").Value = 2196
$ws.Range("K89").Value = 791787.8999999999
$ws.Range("L89").Value = 10980
$ws.Range("M89").Value = -786171.8999999999
$ws.Range("N89").Value = -22212

$ws.Range("H96").Value = 347.05554
$ws.Range("I96").Value = 277.9375
$ws.Range("J96").Value = 900
$ws.Range("K96").Value = 833.8125
$ws.Range("L96").Value = 2700
$ws.Range("M96").Value = 539.1875
$ws.Range("N96").Value = -5446

$ws.Range("H97").Value = 2703.3333
$ws.Range("J97").Value = 2703.3333
$ws.Range("L97").Value = 8109.999899999999
$ws.Range("N97").Value = -9101.999899999999

$ws.Range("H106").Value = 4753
$ws.Range("I106").Value = 4500
$ws.Range("J106").Value = 5006
$ws.Range("K106").Value = 4500
$ws.Range("L106").Value = 5006
$ws.Range("M106").Value = -3869
$ws.Range("N106").Value = -6268

$ws.Range("H129").Value = 1886.75
$ws.Range("I129").Value = 696.75
$ws.Range("K129").Value = 2090.25
$ws.Range("M129").Value = 2909.75

$ws.Range("H135").Value = 956.2143
$ws.Range("I135").Value = 685.7714
$ws.Range("J135").Value = 2308.4285
$ws.Range("K135").Value = 6171.942599999999
$ws.Range("L135").Value = 20775.8565
$ws.Range("M135").Value = -3636.942599999999
$ws.Range("N135").Value = -25845.8565

$ws.Range("H137").Value = 1737.6818
$ws.Range("I137").Value = 3302.923
$ws.Range("J137").Value = 1081.2903
$ws.Range("K137").Value = 9908.769
$ws.Range("L137").Value = 3243.8709
$ws.Range("M137").Value = -7358.769
$ws.Range("N137").Value = -8343.8709

$ws.Range("H138").Value = 2720.1052
$ws.Range("I138").Value = 1694.1034
$ws.Range("J138").Value = 3353.1702
$ws.Range("K138").Value = 5082.3102
$ws.Range("L138").Value = 10059.5106
$ws.Range("M138").Value = 57.6898000000001
$ws.Range("N138").Value = -20339.5106

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3748.33
$ws.Range("I32").Value = 3748.33
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 3748.33
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -3461.33
$ws.Range("N32").ClearContents()

$ws.Range("H88").Value = 946.625
$ws.Range("I88").Value = 898
$ws.Range("J88").Value = 962.8333
$ws.Range("K88").Value = 898
$ws.Range("L88").Value = 962.8333
$ws.Range("M88").Value = -492
$ws.Range("N88").Value = -1774.8333

$ws.Range("H91").Value = 946.625
$ws.Range("I91").Value = 898
$ws.Range("J91").Value = 962.8333
$ws.Range("K91").Value = 898
$ws.Range("L91").Value = 962.8333
$ws.Range("M91").Value = 506
$ws.Range("N91").Value = -3770.8333

$ws.Range("H95").Value = 24537
$ws.Range("J95").Value = 24537
$ws.Range("L95").Value = 24537
$ws.Range("N95").Value = -30029

$ws.Range("H122").Value = 1115.6389
$ws.Range("I122").Value = 974.5833
$ws.Range("J122").Value = 1397.75
$ws.Range("K122").Value = 2923.7499
$ws.Range("L122").Value = 4193.25
$ws.Range("M122").Value = -473.7498999999998
$ws.Range("N122").Value = -9093.25

$ws.Range("H132").Value = 4921.441
$ws.Range("I132").Value = 3716.425
$ws.Range("K132").Value = 11149.275
$ws.Range("M132").Value = -8619.275000000001

$ws.Range("H135").Value = 100020460
$ws.Range("J135").Value = 100020460
$ws.Range("L135").Value = 100020460
$ws.Range("N135").Value = -100030600

$ws.Range("H139").Value = 32878.75
$ws.Range("J139").Value = 32878.75
$ws.Range("L139").Value = 32878.75
$ws.Range("N139").Value = -43158.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1578
$ws.Range("I107").Value = 1596.6923
$ws.Range("J107").Value = 1456.5
$ws.Range("K107").Value = 1596.6923
$ws.Range("L107").Value = 1456.5
$ws.Range("M107").Value = 323.3077000000001
$ws.Range("N107").Value = -5296.5

$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2583.9333
$ws.Range("I31").Value = 1361.1765
$ws.Range("J31").Value = 4182.923
$ws.Range("K31").Value = 1361.1765
$ws.Range("L31").Value = 4182.923
$ws.Range("M31").Value = -1066.1765
$ws.Range("N31").Value = -4772.923

$ws.Range("H34").Value = 2583.9333
$ws.Range("I34").Value = 1361.1765
$ws.Range("J34").Value = 4182.923
$ws.Range("K34").Value = 1361.1765
$ws.Range("L34").Value = 4182.923
$ws.Range("M34").Value = -1159.1765
$ws.Range("N34").Value = -4586.923

$ws.Range("H58").Value = 3148.3257
$ws.Range("I58").Value = 1398.7858
$ws.Range("J58").Value = 6414.1333
$ws.Range("K58").Value = 1398.7858
$ws.Range("L58").Value = 6414.1333
$ws.Range("M58").Value = -1195.7858
$ws.Range("N58").Value = -6820.1333

$ws.Range("H62").Value = 6887.143
$ws.Range("I62").Value = 2079.1667
$ws.Range("J62").Value = 13297.777
$ws.Range("K62").Value = 2079.1667
$ws.Range("L62").Value = 13297.777
$ws.Range("M62").Value = -1455.1667
$ws.Range("N62").Value = -14545.777

$ws.Range("H65").Value = 6887.143
$ws.Range("I65").Value = 2079.1667
$ws.Range("J65").Value = 13297.777
$ws.Range("K65").Value = 10395.8335
$ws.Range("L65").Value = 66488.88499999999
$ws.Range("M65").Value = -7275.833500000001
$ws.Range("N65").Value = -72728.88499999999

$ws.Range("H132").Value = 2239.6553
$ws.Range("I132").Value = 1635.0588
$ws.Range("J132").Value = 3096.1667
$ws.Range("K132").Value = 4905.1764
$ws.Range("L132").Value = 9288.500100000001
$ws.Range("M132").Value = -2375.1764
$ws.Range("N132").Value = -14348.5001

$ws.Range("H134").Value = 1482.8269
$ws.Range("I134").Value = 852.2432
$ws.Range("J134").Value = 3038.2666
$ws.Range("K134").Value = 2556.7296
$ws.Range("L134").Value = 9114.799800000001
$ws.Range("M134").Value = -21.72960000000012
$ws.Range("N134").Value = -14184.7998

$ws.Range("H136").Value = 3148.3257
$ws.Range("I136").Value = 1398.7858
$ws.Range("J136").Value = 6414.1333
$ws.Range("K136").Value = 4196.357400000001
$ws.Range("L136").Value = 19242.3999
$ws.Range("M136").Value = -1646.357400000001
$ws.Range("N136").Value = -24342.3999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 104.210526
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 104.210526
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 312.631578
$ws.Range("M12").ClearContents()
$ws.Range("N12").Value = -658.631578

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2464.818
$ws.Range("I122").Value = 2089.125
$ws.Range("J122").Value = 3466.6667
$ws.Range("K122").Value = 6267.375
$ws.Range("L122").Value = 10400.0001
$ws.Range("M122").Value = -3817.375
$ws.Range("N122").Value = -15300.0001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H95").Value = 21718.77
$ws.Range("J95").Value = 21718.77
$ws.Range("L95").Value = 21718.77
$ws.Range("N95").Value = -27210.77

$ws.Range("H97").Value = 21076.445
$ws.Range("J97").Value = 21076.445
$ws.Range("L97").Value = 21076.445
$ws.Range("N97").Value = -23058.445

$ws.Range("H122").Value = 5387.4116
$ws.Range("I122").Value = 6703.3335
$ws.Range("J122").Value = 3261.6924
$ws.Range("K122").Value = 20110.0005
$ws.Range("L122").Value = 9785.0772
$ws.Range("M122").Value = -17660.0005
$ws.Range("N122").Value = -14685.0772

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H97").Value = 23166.666
$ws.Range("J97").Value = 23166.666
$ws.Range("L97").Value = 23166.666
$ws.Range("N97").Value = -25148.666

$ws.Range("H136").Value = 1613.3098
$ws.Range("I136").Value = 1502.0555
$ws.Range("J136").Value = 1966.7059
$ws.Range("K136").Value = 4506.166499999999
$ws.Range("L136").Value = 5900.1177
$ws.Range("M136").Value = -1956.166499999999
$ws.Range("N136").Value = -11000.1177
